# Weekly update: insert a new Orégano price record for Vega Central
# Mapocho de Santiago ahead of the existing rows (row 42), shifting the
# remaining history down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 42, pushing existing rows 42:84 down to 43:85.
$ws.Rows("42").Insert()

# Populate the new row 42 with this week's data. The descriptive columns
# (market/region/product/unit/etc.) are identical to the rest of the
# series; only the date and price figures change week to week.
$ws.Range("A42").Value2 = 9
$ws.Range("B42").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C42").Value2 = "Metropolitana"
$ws.Range("D42").Value2 = 44895
$ws.Range("E42").Value2 = 13
$ws.Range("F42").Value2 = 100112029
$ws.Range("G42").Value2 = "Orégano"
$ws.Range("H42").Value2 = "Sin especificar"
$ws.Range("I42").Value2 = "Primera"
$ws.Range("J42").Value2 = 16
$ws.Range("K42").Value2 = 17000
$ws.Range("L42").Value2 = 17000
$ws.Range("M42").Value2 = 17000
$ws.Range("N42").Value2 = "$/docena de atados"
$ws.Range("O42").Value2 = "Región Metropolitana"
$ws.Range("P42").Value2 = 5667
$ws.Range("Q42").Value2 = 3
$ws.Range("R42").Value2 = "Hortaliza"
